# edit.ps1 -- applies the betexplorer CFL Group A (Czech Republic 2023-2024) update
# 1) For 8 pairs of adjacent rows, the match-detail columns (F:V) were swapped
#    (index/metadata columns A:E stay with their original row).
# 2) Three new match rows (111-113) were appended at the end of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap match data (columns F:V) between row 2 and row 3 ---
$ws.Range("F2").Value = "Kraluv Dvur"
$ws.Range("F3").Value = "Hostoun"
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("H2").Value = "Domazlice"
$ws.Range("H3").Value = "Ceske Budejovice B"
$ws.Range("I2").Value = 2
$ws.Range("I3").Value = 0
$ws.Range("J2").Value = 3.91
$ws.Range("J3").Value = 1.98
$ws.Range("K2").Value = "03/08/2023 06:12"
$ws.Range("K3").Value = "03/08/2023 06:13"
$ws.Range("L2").Value = 6.38
$ws.Range("L3").Value = 2.04
$ws.Range("M2").Value = "04/08/2023 17:53"
$ws.Range("M3").Value = "04/08/2023 17:05"
$ws.Range("N2").Value = 3.95
$ws.Range("N3").Value = 3.53
$ws.Range("O2").Value = "03/08/2023 06:12"
$ws.Range("O3").Value = "03/08/2023 06:13"
$ws.Range("P2").Value = 4.81
$ws.Range("P3").Value = 3.69
$ws.Range("Q2").Value = "04/08/2023 17:53"
$ws.Range("Q3").Value = "04/08/2023 17:05"
$ws.Range("R2").Value = 1.6
$ws.Range("R3").Value = 2.93
$ws.Range("S2").Value = "03/08/2023 06:12"
$ws.Range("S3").Value = "03/08/2023 06:13"
$ws.Range("T2").Value = 1.39
$ws.Range("T3").Value = 3.08
$ws.Range("U2").Value = "04/08/2023 17:53"
$ws.Range("U3").Value = "04/08/2023 17:05"
$ws.Range("V2").Value = "https://www.betexplorer.com/football/czech-republic/cfl-group-a/kraluv-dvur-domazlice/pKjW3Tsd/"
$ws.Range("V3").Value = "https://www.betexplorer.com/football/czech-republic/cfl-group-a/hostoun-ceske-budejovice/jciS49Sk/"

# --- Swap match data (columns F:V) between row 70 and row 71 ---
$ws.Range("F70").Value = "Dukla Prague B"
$ws.Range("F71").Value = "Bohemians 1905 B"
$ws.Range("G70").Value = 3
$ws.Range("G71").Value = 1
$ws.Range("H70").Value = "Hostoun"
$ws.Range("H71").Value = "Domazlice"
$ws.Range("I70").Value = 3
$ws.Range("I71").Value = 3
$ws.Range("J70").Value = 2.11
$ws.Range("J71").Value = 3.07
$ws.Range("K70").Value = "29/09/2023 21:42"
$ws.Range("K71").Value = "29/09/2023 21:42"
$ws.Range("L70").Value = 2.77
$ws.Range("L71").Value = 4.97
$ws.Range("M70").Value = "01/10/2023 10:12"
$ws.Range("M71").Value = "01/10/2023 10:06"
$ws.Range("N70").Value = 3.56
$ws.Range("N71").Value = 3.59
$ws.Range("O70").Value = "29/09/2023 21:42"
$ws.Range("O71").Value = "29/09/2023 21:42"
$ws.Range("P70").Value = 3.86
$ws.Range("P71").Value = 4.46
$ws.Range("Q70").Value = "01/10/2023 10:12"
$ws.Range("Q71").Value = "01/10/2023 10:06"
$ws.Range("R70").Value = 2.68
$ws.Range("R71").Value = 1.91
$ws.Range("S70").Value = "29/09/2023 21:42"
$ws.Range("S71").Value = "29/09/2023 21:42"
$ws.Range("T70").Value = 2.15
$ws.Range("T71").Value = 1.52
$ws.Range("U70").Value = "01/10/2023 10:12"
$ws.Range("U71").Value = "01/10/2023 10:06"
$ws.Range("V70").Value = "https://www.betexplorer.com/football/czech-republic/cfl-group-a/dukla-prague-hostoun/OS6KlCfr/"
$ws.Range("V71").Value = "https://www.betexplorer.com/football/czech-republic/cfl-group-a/bohemians-1905-domazlice/dEjv9ttc/"

# --- Swap match data (columns F:V) between row 79 and row 80 ---
$ws.Range("F79").Value = "Admira Prague"
$ws.Range("F80").Value = "Slavia Prague B"
$ws.Range("G79").Value = 1
$ws.Range("G80").Value = 1
$ws.Range("H79").Value = "Dukla Prague B"
$ws.Range("H80").Value = "Bohemians 1905 B"
$ws.Range("I79").Value = 1
$ws.Range("I80").Value = 0
$ws.Range("J79").Value = 1.72
$ws.Range("J80").Value = 1.22
$ws.Range("K79").Value = "06/10/2023 21:42"
$ws.Range("K80").Value = "06/10/2023 21:42"
$ws.Range("L79").Value = 1.55
$ws.Range("L80").Value = 1.16
$ws.Range("M79").Value = "08/10/2023 10:05"
$ws.Range("M80").Value = "08/10/2023 09:40"
$ws.Range("N79").Value = 3.89
$ws.Range("N80").Value = 5.98
$ws.Range("O79").Value = "06/10/2023 21:42"
$ws.Range("O80").Value = "06/10/2023 21:42"
$ws.Range("P79").Value = 4.31
$ws.Range("P80").Value = 7.94
$ws.Range("Q79").Value = "08/10/2023 10:05"
$ws.Range("Q80").Value = "08/10/2023 09:51"
$ws.Range("R79").Value = 3.45
$ws.Range("R80").Value = 6.85
$ws.Range("S79").Value = "06/10/2023 21:42"
$ws.Range("S80").Value = "06/10/2023 21:42"
$ws.Range("T79").Value = 4.81
$ws.Range("T80").Value = 10.16
$ws.Range("U79").Value = "08/10/2023 10:05"
$ws.Range("U80").Value = "08/10/2023 09:51"
$ws.Range("V79").Value = "https://www.betexplorer.com/football/czech-republic/cfl-group-a/admira-prague-dukla-prague/rm4Q7i1K/"
$ws.Range("V80").Value = "https://www.betexplorer.com/football/czech-republic/cfl-group-a/slavia-prague-bohemians-1905/0ObV6BGQ/"

# --- Swap match data (columns F:V) between row 83 and row 84 ---
$ws.Range("F83").Value = "Motorlet Prague"
$ws.Range("F84").Value = "Bohemians 1905 B"
$ws.Range("G83").Value = 2
$ws.Range("G84").Value = 0
$ws.Range("H83").Value = "Admira Prague"
$ws.Range("H84").Value = "Hostoun"
$ws.Range("I83").Value = 1
$ws.Range("I84").Value = 1
$ws.Range("J83").Value = 2.35
$ws.Range("J84").Value = 1.96
$ws.Range("K83").Value = "12/10/2023 21:42"
$ws.Range("K84").Value = "12/10/2023 21:42"
$ws.Range("L83").Value = 2.64
$ws.Range("L84").Value = 2.39
$ws.Range("M83").Value = "14/10/2023 10:13"
$ws.Range("M84").Value = "14/10/2023 10:10"
$ws.Range("N83").Value = 3.33
$ws.Range("N84").Value = 3.63
$ws.Range("O83").Value = "12/10/2023 21:42"
$ws.Range("O84").Value = "12/10/2023 21:42"
$ws.Range("P83").Value = 3.6
$ws.Range("P84").Value = 3.81
$ws.Range("Q83").Value = "14/10/2023 09:43"
$ws.Range("Q84").Value = "14/10/2023 10:10"
$ws.Range("R83").Value = 2.49
$ws.Range("R84").Value = 2.92
$ws.Range("S83").Value = "12/10/2023 21:42"
$ws.Range("S84").Value = "12/10/2023 21:42"
$ws.Range("T83").Value = 2.34
$ws.Range("T84").Value = 2.48
$ws.Range("U83").Value = "14/10/2023 10:13"
$ws.Range("U84").Value = "14/10/2023 10:10"
$ws.Range("V83").Value = "https://www.betexplorer.com/football/czech-republic/cfl-group-a/motorlet-prague-admira-prague/dfp11kvl/"
$ws.Range("V84").Value = "https://www.betexplorer.com/football/czech-republic/cfl-group-a/bohemians-1905-hostoun/WCxEbmO6/"

# --- Swap match data (columns F:V) between row 94 and row 95 ---
$ws.Range("F94").Value = "Admira Prague"
$ws.Range("F95").Value = "Slavia Prague B"
$ws.Range("G94").Value = 0
$ws.Range("G95").Value = 3
$ws.Range("H94").Value = "Bohemians 1905 B"
$ws.Range("H95").Value = "Taborsko akademie"
$ws.Range("I94").Value = 0
$ws.Range("I95").Value = 1
$ws.Range("J94").Value = 2.39
$ws.Range("J95").Value = 1.36
$ws.Range("K94").Value = "22/10/2023 08:27"
$ws.Range("K95").Value = "22/10/2023 08:27"
$ws.Range("L94").Value = 2.11
$ws.Range("L95").Value = 1.31
$ws.Range("M94").Value = "22/10/2023 10:00"
$ws.Range("M95").Value = "22/10/2023 10:01"
$ws.Range("N94").Value = 3.45
$ws.Range("N95").Value = 5.36
$ws.Range("O94").Value = "22/10/2023 08:27"
$ws.Range("O95").Value = "22/10/2023 08:27"
$ws.Range("P94").Value = 3.5
$ws.Range("P95").Value = 5.53
$ws.Range("Q94").Value = "22/10/2023 10:00"
$ws.Range("Q95").Value = "22/10/2023 10:01"
$ws.Range("R94").Value = 2.66
$ws.Range("R95").Value = 6.11
$ws.Range("S94").Value = "22/10/2023 08:27"
$ws.Range("S95").Value = "22/10/2023 08:27"
$ws.Range("T94").Value = 3.07
$ws.Range("T95").Value = 6.75
$ws.Range("U94").Value = "22/10/2023 10:00"
$ws.Range("U95").Value = "22/10/2023 10:01"
$ws.Range("V94").Value = "https://www.betexplorer.com/football/czech-republic/cfl-group-a/admira-prague-bohemians-1905/WW0eBQ7C/"
$ws.Range("V95").Value = "https://www.betexplorer.com/football/czech-republic/cfl-group-a/slavia-prague-taborsko-akademie/Eej0ApNI/"

# --- Swap match data (columns F:V) between row 98 and row 99 ---
$ws.Range("F98").Value = "Pisek"
$ws.Range("F99").Value = "Motorlet Prague"
$ws.Range("G98").Value = 1
$ws.Range("G99").Value = 2
$ws.Range("H98").Value = "Admira Prague"
$ws.Range("H99").Value = "Vltavin"
$ws.Range("I98").Value = 1
$ws.Range("I99").Value = 2
$ws.Range("J98").Value = 1.79
$ws.Range("J99").Value = 2.2
$ws.Range("K98").Value = "27/10/2023 22:13"
$ws.Range("K99").Value = "27/10/2023 22:13"
$ws.Range("L98").Value = 2.63
$ws.Range("L99").Value = 2.09
$ws.Range("M98").Value = "28/10/2023 10:04"
$ws.Range("M99").Value = "28/10/2023 09:54"
$ws.Range("N98").Value = 3.78
$ws.Range("N99").Value = 3.48
$ws.Range("O98").Value = "27/10/2023 22:13"
$ws.Range("O99").Value = "27/10/2023 22:13"
$ws.Range("P98").Value = 3.78
$ws.Range("P99").Value = 3.48
$ws.Range("Q98").Value = "28/10/2023 10:04"
$ws.Range("Q99").Value = "28/10/2023 09:54"
$ws.Range("R98").Value = 3.58
$ws.Range("R99").Value = 2.79
$ws.Range("S98").Value = "27/10/2023 22:13"
$ws.Range("S99").Value = "27/10/2023 22:13"
$ws.Range("T98").Value = 2.27
$ws.Range("T99").Value = 3.12
$ws.Range("U98").Value = "28/10/2023 10:04"
$ws.Range("U99").Value = "28/10/2023 09:54"
$ws.Range("V98").Value = "https://www.betexplorer.com/football/czech-republic/cfl-group-a/pisek-admira-prague/hnk494xP/"
$ws.Range("V99").Value = "https://www.betexplorer.com/football/czech-republic/cfl-group-a/motorlet-prague-loko-vltavin/t2rnYmFP/"

# --- Swap match data (columns F:V) between row 102 and row 103 ---
$ws.Range("F102").Value = "Dukla Prague B"
$ws.Range("F103").Value = "Bohemians 1905 B"
$ws.Range("G102").Value = 3
$ws.Range("G103").Value = 0
$ws.Range("H102").Value = "Povltavska FA"
$ws.Range("H103").Value = "Kraluv Dvur"
$ws.Range("I102").Value = 3
$ws.Range("I103").Value = 0
$ws.Range("J102").Value = 2.29
$ws.Range("J103").Value = 1.5
$ws.Range("K102").Value = "29/10/2023 01:42"
$ws.Range("K103").Value = "29/10/2023 01:42"
$ws.Range("L102").Value = 2.15
$ws.Range("L103").Value = 1.7
$ws.Range("M102").Value = "29/10/2023 10:14"
$ws.Range("M103").Value = "29/10/2023 10:07"
$ws.Range("N102").Value = 3.49
$ws.Range("N103").Value = 4.57
$ws.Range("O102").Value = "29/10/2023 01:42"
$ws.Range("O103").Value = "29/10/2023 01:42"
$ws.Range("P102").Value = 3.61
$ws.Range("P103").Value = 4.07
$ws.Range("Q102").Value = "29/10/2023 10:05"
$ws.Range("Q103").Value = "29/10/2023 10:07"
$ws.Range("R102").Value = 2.65
$ws.Range("R103").Value = 4.59
$ws.Range("S102").Value = "29/10/2023 01:42"
$ws.Range("S103").Value = "29/10/2023 01:42"
$ws.Range("T102").Value = 2.91
$ws.Range("T103").Value = 3.96
$ws.Range("U102").Value = "29/10/2023 10:14"
$ws.Range("U103").Value = "29/10/2023 10:06"
$ws.Range("V102").Value = "https://www.betexplorer.com/football/czech-republic/cfl-group-a/dukla-prague-povltavska-fa/vDhROqxC/"
$ws.Range("V103").Value = "https://www.betexplorer.com/football/czech-republic/cfl-group-a/bohemians-1905-kraluv-dvur/YaNQQ57a/"

# --- Swap match data (columns F:V) between row 108 and row 109 ---
$ws.Range("F108").Value = "Ceske Budejovice B"
$ws.Range("F109").Value = "Kraluv Dvur"
$ws.Range("G108").Value = 0
$ws.Range("G109").Value = 1
$ws.Range("H108").Value = "Slavia Prague B"
$ws.Range("H109").Value = "Pisek"
$ws.Range("I108").Value = 0
$ws.Range("I109").Value = 5
$ws.Range("J108").Value = 2.89
$ws.Range("J109").Value = 1.97
$ws.Range("K108").Value = "04/11/2023 01:13"
$ws.Range("K109").Value = "04/11/2023 01:13"
$ws.Range("L108").Value = 3.28
$ws.Range("L109").Value = 1.98
$ws.Range("M108").Value = "04/11/2023 09:37"
$ws.Range("M109").Value = "04/11/2023 10:27"
$ws.Range("N108").Value = 3.8
$ws.Range("N109").Value = 3.59
$ws.Range("O108").Value = "04/11/2023 01:13"
$ws.Range("O109").Value = "04/11/2023 01:13"
$ws.Range("P108").Value = 4.13
$ws.Range("P109").Value = 3.53
$ws.Range("Q108").Value = "04/11/2023 10:25"
$ws.Range("Q109").Value = "04/11/2023 10:27"
$ws.Range("R108").Value = 2.03
$ws.Range("R109").Value = 3.16
$ws.Range("S108").Value = "04/11/2023 01:13"
$ws.Range("S109").Value = "04/11/2023 01:13"
$ws.Range("T108").Value = 1.85
$ws.Range("T109").Value = 3.37
$ws.Range("U108").Value = "04/11/2023 10:25"
$ws.Range("U109").Value = "04/11/2023 10:27"
$ws.Range("V108").Value = "https://www.betexplorer.com/football/czech-republic/cfl-group-a/ceske-budejovice-slavia-prague/WA7oQgYo/"
$ws.Range("V109").Value = "https://www.betexplorer.com/football/czech-republic/cfl-group-a/kraluv-dvur-pisek/CQjZMN6O/"

# --- Append new rows 111-113 ---
# Row 111
$ws.Range("A109").Copy()
$ws.Range("A111").PasteSpecial(-4122)
$ws.Range("E109").Copy()
$ws.Range("E111").PasteSpecial(-4122)
$ws.Range("A111").Value = 110
$ws.Range("B111").Value = "czech-republic"
$ws.Range("C111").Value = "cfl-group-a"
$ws.Range("D111").Value = "2023-2024"
$ws.Range("E111").Value = 45235.42708333334
$ws.Range("F111").Value = "Admira Prague"
$ws.Range("G111").Value = 3
$ws.Range("H111").Value = "Taborsko akademie"
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 1.65
$ws.Range("K111").Value = "05/11/2023 08:41"
$ws.Range("L111").Value = 1.65
$ws.Range("M111").Value = "05/11/2023 08:41"
$ws.Range("N111").Value = 3.95
$ws.Range("O111").Value = "05/11/2023 08:41"
$ws.Range("P111").Value = 3.95
$ws.Range("Q111").Value = "05/11/2023 08:41"
$ws.Range("R111").Value = 4.39
$ws.Range("S111").Value = "05/11/2023 08:41"
$ws.Range("T111").Value = 4.39
$ws.Range("U111").Value = "05/11/2023 08:41"
$ws.Range("V111").Value = "https://www.betexplorer.com/football/czech-republic/cfl-group-a/admira-prague-taborsko-akademie/Ma0bNiI4/"

# Row 112
$ws.Range("A109").Copy()
$ws.Range("A112").PasteSpecial(-4122)
$ws.Range("E109").Copy()
$ws.Range("E112").PasteSpecial(-4122)
$ws.Range("A112").Value = 111
$ws.Range("B112").Value = "czech-republic"
$ws.Range("C112").Value = "cfl-group-a"
$ws.Range("D112").Value = "2023-2024"
$ws.Range("E112").Value = 45235.4375
$ws.Range("F112").Value = "Plzen B"
$ws.Range("G112").Value = 2
$ws.Range("H112").Value = "Karlovy Vary"
$ws.Range("I112").Value = 1
$ws.Range("J112").Value = 1.33
$ws.Range("K112").Value = "05/11/2023 07:41"
$ws.Range("L112").Value = 1.33
$ws.Range("M112").Value = "05/11/2023 07:41"
$ws.Range("N112").Value = 5.13
$ws.Range("O112").Value = "05/11/2023 08:30"
$ws.Range("P112").Value = 5.13
$ws.Range("Q112").Value = "05/11/2023 08:30"
$ws.Range("R112").Value = 7.13
$ws.Range("S112").Value = "05/11/2023 07:41"
$ws.Range("T112").Value = 7.13
$ws.Range("U112").Value = "05/11/2023 07:41"
$ws.Range("V112").Value = "https://www.betexplorer.com/football/czech-republic/cfl-group-a/plzen-karlovy-vary/vyBgOX2b/"

# Row 113
$ws.Range("A109").Copy()
$ws.Range("A113").PasteSpecial(-4122)
$ws.Range("E109").Copy()
$ws.Range("E113").PasteSpecial(-4122)
$ws.Range("A113").Value = 112
$ws.Range("B113").Value = "czech-republic"
$ws.Range("C113").Value = "cfl-group-a"
$ws.Range("D113").Value = "2023-2024"
$ws.Range("E113").Value = 45235.58333333334
$ws.Range("F113").Value = "Povltavska FA"
$ws.Range("G113").Value = 1
$ws.Range("H113").Value = "Motorlet Prague"
$ws.Range("I113").Value = 2
$ws.Range("J113").Value = 1.47
$ws.Range("K113").Value = "05/11/2023 12:03"
$ws.Range("L113").Value = 1.49
$ws.Range("M113").Value = "05/11/2023 13:57"
$ws.Range("N113").Value = 4.59
$ws.Range("O113").Value = "05/11/2023 12:03"
$ws.Range("P113").Value = 4.29
$ws.Range("Q113").Value = "05/11/2023 13:57"
$ws.Range("R113").Value = 4.98
$ws.Range("S113").Value = "05/11/2023 12:03"
$ws.Range("T113").Value = 5.53
$ws.Range("U113").Value = "05/11/2023 13:57"
$ws.Range("V113").Value = "https://www.betexplorer.com/football/czech-republic/cfl-group-a/povltavska-fa-motorlet-prague/0pCkPDmh/"

